# 09_wa_automation/contact.xlsx — trim the contact list down to a single
# test row and point it at a generic placeholder number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the third contact (Rupesh / 7857831014) entirely — the sheet now
# only tracks one row of sample data.
$ws.Rows.Item(3).Delete()

# Replace the remaining contact's number with a generic placeholder and
# shorten the display name from "Rupesh So" to "Rupesh".
$ws.Range("A2").Value = 1234567890
$ws.Range("B2").Value = "Rupesh"

# Column A widens slightly now that "bestFit" recalculates for the
# remaining content.
$ws.Columns.Item(1).ColumnWidth = 11.166666666666666

# Leave the selection on D9, matching the saved view state.
$ws.Range("D9").Select()
